$d = $word.ActiveDocument

# 1. Strike through the "Tie empty cell generation to number of bees" bullet.
$pTie = $d.Paragraphs(4)
if ($pTie.Range.Text -notmatch "Tie empty cell generation") {
    throw "Paragraph 4 text mismatch: $($pTie.Range.Text)"
}
$pTie.Range.Font.StrikeThrough = 1

# 2. Remove the (stray) strike-through from the "Bees keep hive clean..." bullet,
#    leaving the paragraph mark / other formatting untouched.
$pBees = $d.Paragraphs(6)
if ($pBees.Range.Text -notmatch "Bees keep hive clean") {
    throw "Paragraph 6 text mismatch: $($pBees.Range.Text)"
}
$rBees = $pBees.Range
$rBees.End = $rBees.End - 1
$rBees.Font.StrikeThrough = 0

# 3. Update the "Queens" bullet text.
$d.Content.Find.Execute("Queens - mated (y/m), burging (y/n)", $false, $true, $false, $false, $false, $true, 1, $false, "Queens - mated (y/n)\", 2)

# 4. Insert a new bullet ("Graphical tracking of values (history)") right after
#    the Queens bullet, matching its list/paragraph formatting.
$pQueens = $d.Paragraphs(7)
if ($pQueens.Range.Text -notmatch "Queens - mated") {
    throw "Paragraph 7 text mismatch: $($pQueens.Range.Text)"
}
$pQueens.Range.InsertParagraphAfter()
$pNew = $d.Paragraphs(8)
$pNew.Range.Text = "Graphical tracking of values (history)"
